$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper: apply center/center alignment (+ optional number format) to a range
# ---------------------------------------------------------------------------
function Set-CenterFormat($range, $numFmt) {
    if ($numFmt) { $range.NumberFormat = $numFmt }
    $range.HorizontalAlignment = -4108
    $range.VerticalAlignment = -4108
}

# ---------------------------------------------------------------------------
# 1. Re-format the existing "Best Case" / "Worst Case" block (rows 2-9,
#    columns B-H) so every populated cell is center+center aligned. The
#    decimal-time columns (D and G) additionally keep the 0.000 number
#    format.
# ---------------------------------------------------------------------------
Set-CenterFormat $ws.Range("B2:H3") $null
Set-CenterFormat $ws.Range("B4:C9") $null
Set-CenterFormat $ws.Range("E4:F8") $null
Set-CenterFormat $ws.Range("E9") $null
Set-CenterFormat $ws.Range("H4:H8") $null
Set-CenterFormat $ws.Range("D3:D9") "0.000"
Set-CenterFormat $ws.Range("G3:G8") "0.000"
Set-CenterFormat $ws.Range("B10") $null

# New data point that was previously missing.
$ws.Range("G8").Value = 1149.732

# ---------------------------------------------------------------------------
# 2. New "Difference" / "Growth" blocks in columns I:L
# ---------------------------------------------------------------------------
$ws.Range("I2").Value = "Difference"
$ws.Range("I2:J2").Merge()
$ws.Range("K2").Value = "Growth"
$ws.Range("K2:L2").Merge()
Set-CenterFormat $ws.Range("I2:L2") $null

$ws.Range("I3").Value = "Time"
$ws.Range("J3").Value = "Positions"
$ws.Range("K3").Value = "Time"
$ws.Range("L3").Value = "Positions"
Set-CenterFormat $ws.Range("I3") $null
Set-CenterFormat $ws.Range("J3") "0"
Set-CenterFormat $ws.Range("K3:L3") $null

# Difference formulas (Worst Case - Best Case)
$ws.Range("I4").Formula = "=G4-D4"
$ws.Range("J4").Formula = "=H4-E4"
Set-CenterFormat $ws.Range("I4") $null
Set-CenterFormat $ws.Range("J4") "0"

for ($r = 5; $r -le 9; $r++) {
    $ws.Range("I$r").Formula = "=G$r-D$r"
    $ws.Range("J$r").Formula = "=H$r-E$r"
    Set-CenterFormat $ws.Range("I$r") $null
    Set-CenterFormat $ws.Range("J$r") "0"
}

# Growth formulas (ratio vs previous depth). Formula entry can inherit the
# number format of its precedents (mirrors Excel's own "format came along
# for the ride" behaviour) - reset to the base "Normal" style first so the
# result stays General instead of picking up 0.000 from the I/J columns.
for ($r = 5; $r -le 9; $r++) {
    $prev = $r - 1
    $ws.Range("K$r").Formula = "=I$r/I$prev"
    $ws.Range("L$r").Formula = "=J$r/J$prev"
    $ws.Range("K$r").Style = "Normal"
    $ws.Range("L$r").Style = "Normal"
    Set-CenterFormat $ws.Range("K$r") $null
    Set-CenterFormat $ws.Range("L$r") $null
}

# ---------------------------------------------------------------------------
# 3. Remove the trailing depth-only rows (8-15), which are no longer needed
# ---------------------------------------------------------------------------
$ws.Range("B11:B18").EntireRow.Delete()

# ---------------------------------------------------------------------------
# 4. Restore the view state
# ---------------------------------------------------------------------------
$ws.Range("G15").Select()
